# Split the single "Right click on Grader.java (...)" run into several
# runs separated by w:proofErr markers (spell/grammar-check boundaries),
# matching the sentence-formation cleanup described in the commit
# message, without altering the visible text.

$d = $word.ActiveDocument

$oldText = "Right click on Grader.java (CS421-ProjectP1/src/(default package)"

# Locate the exact run of text to replace using Find (keeps us independent
# of hard-coded character offsets).
$target = $d.Content
$target.Find.ClearFormatting()
$found = $target.Find.Execute($oldText, $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target sentence in document"
}

# Collapse the matched range to an insertion point, removing the old text.
$target.Text = ""

# New content: the same text, re-flowed across multiple runs with
# w:proofErr spell/grammar boundary markers in between.
$newBody = (
    '<w:r><w:t>Right click on Grader.java (CS421-ProjectP1/</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>src</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>/(</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>default package)</w:t></w:r>'
)

$paraOpen = '<w:p w:rsidR="008E500D" w:rsidRDefault="008E500D" w:rsidP="00833B7B">'

$xml = (
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" ' +
    'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<pkg:part pkg:name="/word/document.xml" ' +
    'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $paraOpen + $newBody + '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData>' +
    '</pkg:part>' +
    '</pkg:package>'
)

$target.InsertXML($xml)
